$wb = $excel.ActiveWorkbook

# "Repayment Schedule" is the 3rd sheet (Input, Summary, Repayment Schedule, Transactions)
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting N:P -> O:Q.
# This mirrors a normal Excel "Insert Column" action, which copies the
# formatting of the column immediately to its left (column M).
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab and update its selection.
$ws.Activate()
$ws.Range("S13").Select()
